$wb = $excel.ActiveWorkbook

$mdTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f90fd891fe89046f7c575c6a21930debdef12886/e2e/60fbb8d4-661e-4a21-9688-2ae0bcc04d4b.md"
$mdDisplay = "60fbb8d4-661e-4a21-9688-2ae0bcc04d4b.md"

# Overview sheet: Status (handback) text shown in zh-cn / de-de columns
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "Handed back: in sync with en-US"
$ws1.Range("F2").Value = "Handed back: in sync with en-US"

# zh-cn sheet
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws2.Range("I2").Value = $mdDisplay
$ws2.Range("J2").Value = "60fbb8d4-661e-4a21-9688-2ae0bcc04d4b.e9ffd5825af93b237f144eebeaf01e3b05938745.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-08-16 00:54:00"
$ws2.Hyperlinks.Add($ws2.Range("I2"), $mdTarget, "", "", $mdDisplay)

# de-de sheet
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("I2").Value = $mdDisplay
$ws3.Range("J2").Value = "60fbb8d4-661e-4a21-9688-2ae0bcc04d4b.e9ffd5825af93b237f144eebeaf01e3b05938745.de-de.xlf"
$ws3.Range("K2").Value = "2016-08-16 00:54:12"
$ws3.Hyperlinks.Add($ws3.Range("I2"), $mdTarget, "", "", $mdDisplay)
